# fix: Stringify trust_name and rights_type aggregations.
#
# Columns D (trust_name) and E (rights_type) hold Python/numpy-style
# stringified list reprs, e.g. "['Common Schools' 'Dept. of Transportation']".
# Convert each to a plain comma-separated string:
#   "Common Schools, Dept. of Transportation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

function Convert-ListRepr($raw) {
    if ($raw -eq $null) { return $raw }
    if ($raw.GetType().Name -ne "String") { return $raw }
    if ($raw.Length -eq 0) { return $raw }
    if ($raw.Substring(0,1) -ne "[") { return $raw }

    $itemMatches = [regex]::Matches($raw, "'(.*?)'")
    if ($itemMatches.Count -eq 0) { return $raw }

    $items = @()
    foreach ($m in $itemMatches) {
        $val = $m.Groups[1].Value
        $val = $val -replace "\s+", " "
        $val = $val.Trim()
        $items += $val
    }

    return ($items -join ", ")
}

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)

    $dVal = $dCell.Value()
    $eVal = $eCell.Value()

    $dNew = Convert-ListRepr $dVal
    $eNew = Convert-ListRepr $eVal

    if ($dNew -ne $dVal) { $dCell.Value = $dNew }
    if ($eNew -ne $eVal) { $eCell.Value = $eNew }
}
